$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before B (this shifts B..V -> C..W, copying styles/number
#    formats from the column to the left, and Excel auto-adjusts formulas).
$ws.Columns("B").Insert()

# 2. Remove the old "# paq" placeholder column, which after the insert above now
#    sits at column S (between "Precio" and "Subtotal"). Deleting it shifts
#    T..W back to S..V, restoring the original A1:V10 dimension.
$ws.Columns("S").Delete()

# 3. New column B is "FechaVencimiento" - a due-date column next to "FECHA".
#    Header: bold + centered.
$ws.Range("B1").Value = "FechaVencimiento"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108  # xlCenter

# Populate the due-date values on the three data rows (mirrors column A's date).
$ws.Range("B3").Value = $ws.Range("A3").Value()
$ws.Range("B7").Value = $ws.Range("A7").Value()
$ws.Range("B10").Value = $ws.Range("A10").Value()

# 4. Rename a few headers.
$ws.Range("K1").Value = "Guia"

$ws.Range("L1").Value = "NumeroPaquetes"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").HorizontalAlignment = -4108  # xlCenter

$ws.Range("M1").Value = "Empaco"
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").HorizontalAlignment = -4108  # xlCenter

# 5. Update client reference numbers (G column) for the three data rows.
$ws.Range("G3").Value = 8104442
$ws.Range("G7").Value = 8104442
$ws.Range("G10").Value = 8104441

# 6. Fix up the Subtotal formulas so they reference the now-shifted Cantidad/Precio
#    columns (Q*R instead of the old P*Q); Excel already adjusted these
#    automatically on insert/delete, but set them explicitly to be safe.
$ws.Range("S3").Formula = "=Q3*R3"
$ws.Range("S4").Formula = "=Q4*R4"
$ws.Range("S5").Formula = "=Q5*R5"
$ws.Range("S7").Formula = "=Q7*R7"
$ws.Range("S8").Formula = "=Q8*R8"
$ws.Range("S10").Formula = "=Q10*R10"

$ws.Range("R7").Formula = "=IF(F7=`"`",U7, ROUND(U7/(1+F7),0))"
$ws.Range("R8").Formula = "=IF(F8=`"`",U8, ROUND(U8/(1+F8),0))"
$ws.Range("R10").Formula = "=IF(F10=`"`",U10, ROUND(U10/(1+F10),0))"

# 7. Update the active selection to match the saved state.
$ws.Range("G3").Select()
